# Rewrites the "A ser preenchido" template sheet so each column header is a
# struct attribute (Name, Age, Email, Games) and the dropdown validation that
# used to live on the Month/Year columns now lives on the new Games column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old Month/Year list validations (B2:B1000, C2:C1000).
$ws.Range("B2:B1000").Validation.Delete()
$ws.Range("C2:C1000").Validation.Delete()

# Header row now reflects the struct fields.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Games"

# New dropdown validation for the Games column.
$gamesRange = $ws.Range("D2:D100000")
$gamesRange.Validation.Add(3, 1, 1, '"Super Mario,SONIC,Zelda,GTA"')
$gamesRange.Validation.ShowInput = $false
$gamesRange.Validation.ShowError = $false
